$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the old "July" column (I) / old Total column, collapsing the
# "Total" header + formulas into column H (SUM of Jan:Jun).
$ws.Range("I1:I7").Delete() | Out-Null

# Header row
$ws.Range("H1").Value = "Total"

# New monthly figures for "2025 Gross Profit" (row 2)
$ws.Range("B2").Value = 124559.84
$ws.Range("C2").Value = 116465.91
$ws.Range("D2").Value = 102788.45
$ws.Range("E2").Value = 124082.6
$ws.Range("F2").Value = 101734.84
$ws.Range("G2").Value = 104787.57
$ws.Range("H2").Formula = "=SUM(B2:G2)"

# New monthly figures for "2025 Total Expenses" (row 3)
$ws.Range("B3").Value = 122998.31
$ws.Range("C3").Value = 127047.08
$ws.Range("D3").Value = 124268.94
$ws.Range("E3").Value = 131430.52
$ws.Range("F3").Value = 156847.55
$ws.Range("G3").Value = 116550.34
$ws.Range("H3").Formula = "=SUM(B3:G3)"

# Row 4 "2025 Margin" shared formula now spans only B4:H4
$ws.Range("B4").Formula = "=(B2-B3)/B2"
$ws.Range("C4:H4").Formula = "=(C2-C3)/C2"

# Row 5 "2024 Margin" shared formula now spans only B5:H5
$ws.Range("B5").Formula = "=(B6-B7)/B6"
$ws.Range("C5:H5").Formula = "=(C6-C7)/C6"

# New monthly figures for "2024 Gross Profit" (row 6)
$ws.Range("B6").Value = 89259.06
$ws.Range("C6").Value = 77299.86
$ws.Range("D6").Value = 72777.96
$ws.Range("E6").Value = 66949.03
$ws.Range("F6").Value = 77545.02
$ws.Range("G6").Value = 70120.57
$ws.Range("H6").Formula = "=SUM(B6:G6)"

# New monthly figures for "2024 Total Expenses" (row 7)
$ws.Range("B7").Value = 118497.98
$ws.Range("C7").Value = 108808.14
$ws.Range("D7").Value = 109359.16
$ws.Range("E7").Value = 138064.01
$ws.Range("F7").Value = 136170.16
$ws.Range("G7").Value = 118492.83
$ws.Range("H7").Formula = "=SUM(B7:G7)"

# Row 2 no longer has the taller custom row height - restore autofit default
$ws.Rows(2).AutoFit() | Out-Null
